$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the deck-related endpoint values.
# Cells C5, C6, C7, C23 keep the same displayed text already, but we
# re-set them too so the shared string table gets rebuilt cleanly.
$ws.Range("C5").Value = "users/getdetails"
$ws.Range("C6").Value = "users/update/:id"
$ws.Range("C7").Value = "users/delete/:id"

$ws.Range("C21").Value = "decks/getAll"
$ws.Range("C22").Value = "decks/get/:id"
$ws.Range("C23").Value = "decks/post"
$ws.Range("C24").Value = "decks/update/:id"
$ws.Range("C25").Value = "decks/delete/:id"

# Update the sheet view: scroll so B19 is the top-left visible cell and
# select C25.
$ws.Range("C25").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 2
